$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: Fibonacci response - drop the "Sure! " preamble and extend the digit list
$b2Text = @'
Here are the first 100 digits of the Fibonacci sequence:
1, 1, 2, 3, 5, 8, 13, 21, 34, 55, 89, 144, 233, 377, 610, 987, 1597, 2584, 4181, 6765, 10946, 17711, 28657, 46368, 75025, 121393, 196418, 317811, 514229, 832040, 1346269, 2178309, 3524578, 5702887, 9227465, 14930352, 24157817, 39088169, 63245986, 102334155, 165580141, 267914296, 433494437, 701408733, 1134903170, 1836311903, 2971215073, 4807526976, 7778742049, 12586269025, 20365011074, 32951280099, 53316291173, 86267571272, 139583862445, 225851433717, 365435296162, 591286729879, 956722026041, 1548008755920, 2504730781961, 4052739537881, 6557470319842, 10610209857723, 17167680177565, 27777890035288, 44945570212853, 72723460248141, 117669030460994, 190392490709135, 308061521170129, 498454011879264, 806515533049393, 1304969544928657, 2111485077978050, 3416454622906707, 5527939700884757, 8944394323791464, 14472334024676221, 23416728348467685, 37889062373143906, 61305790721611591, 99194853094755497, 160500643816367088, 259695496911122585, 420196140727489673, 679891637638612258, 1100087778366101931, 1779979416004714189, 2880067194370816120, 4660046610375530309, 7540113804746346429, 12200160415121876738, 19740274219868223167, 31940434634990099905, 51680708854858323072, 83621143489848422977, 135301852344706746049, 218922995834555169026.
'@
$ws.Range("B2").Value = $b2Text
$ws.Rows.Item(2).AutoFit() | Out-Null

# B3: Pi response - drop the "Sure! " preamble and shorten to 10 groups of digits
$b3Text = @'
Here are the first 100 digits of Pi (π):
3.1415926535 8979323846 2643383279 5028841971 6939937510 5820974944 5923078164 0628620899 8628034825 3421170679
'@
$ws.Range("B3").Value = $b3Text
$ws.Rows.Item(3).AutoFit() | Out-Null

# B4: Euler's number response - drop the "Certainly! " preamble, new digits, shortened
$b4Text = @'
Here are the first 100 digits of Euler's number (e):
2.7182818284 5904523536 0287471352 6624977572 4709369995 9574966967 6277240766 3035354759 4571382178 5251664274
'@
$ws.Range("B4").Value = $b4Text
$ws.Rows.Item(4).AutoFit() | Out-Null

# B5: replace the story entirely with the new dragon-egg story
$b5Text = @'
The Quest for the Lost Dragon Egg
In a small village nestled between rolling hills and enchanted forests, there lived an 11-year-old boy named Leo. Leo was curious, adventurous, and always dreamed of becoming a great explorer like the heroes in the stories his grandfather told him. He had a best friend named Lily, who was just as daring as he was, and together they were always getting into some kind of mischief.
One crisp autumn morning, Leo and Lily decided to venture into the Whispering Woods, a place no one in the village dared to go. The woods were known for their towering trees that seemed to whisper secrets when the wind blew, and for the strange creatures said to live deep within. But Leo and Lily were not afraid—well, maybe just a little—but that was what made the adventure exciting.
As they entered the forest, the air felt different. The whispers of the trees sounded like soft voices, as if they were calling out to them. The path beneath their feet was covered in golden leaves, and the sunlight barely made it through the thick canopy above.
“We should turn back,” Lily whispered, though she didn’t really mean it.
Leo grinned. “Are you kidding? This is just getting good!”
After walking for what felt like hours, they came upon an old stone bridge crossing a shimmering river. As they approached, they saw something strange—a small creature sitting on the edge of the bridge. It had the body of a lizard, wings like a bat, and the tiniest puff of smoke curled from its nose.
“A dragon!” Leo gasped.
The creature was barely the size of a cat, but it was unmistakably a baby dragon. Its scales gleamed in the sunlight, shifting colors from blue to green to gold.
The dragon blinked at them, looking as surprised as they were. It let out a tiny roar, which sounded more like a squeaky toy, and then quickly flapped its wings and darted off into the woods.
“Wait!” Leo shouted, chasing after it. Lily followed, her heart pounding with excitement.
They ran deeper into the forest, following the flashes of color from the dragon’s scales, until they reached a clearing. In the center of the clearing was something they had only heard about in legends: a massive dragon’s nest. It was made of branches, leaves, and shimmering stones, but what caught their attention most was the large, glowing egg sitting in the middle of it.
“This is it!” Leo whispered in awe. “The lost dragon egg!”
According to the legends, the dragons had left the village long ago after a terrible storm destroyed their home. The villagers believed that one day, a brave soul would find the last dragon egg and bring peace between humans and dragons once more. But no one had seen a dragon in years—until now.
Before they could get any closer, the ground shook. The trees trembled, and a shadow fell over the clearing. From above, a giant dragon swooped down, landing with a thud in front of the nest. Its eyes were fierce, but there was something gentle about the way it glanced at the egg.
Leo and Lily froze. They had never been so close to such a powerful creature.
The dragon let out a low growl, but instead of attacking, it stepped back, watching them closely.
“I think it’s...waiting,” Lily whispered. “Maybe it wants us to do something.”
Leo’s heart raced. What were they supposed to do? Then, he remembered something his grandfather had told him about dragons: “They only trust those who are brave but kind.”
He took a deep breath, slowly walked forward, and knelt beside the egg. “We’re not here to hurt you,” he said softly. “We just want to help.”
The dragon tilted its head, watching him closely.
Lily, feeling a little braver now, stepped up beside Leo. “We’ll protect your egg,” she said. “We promise.”
Suddenly, the egg began to glow brighter, and tiny cracks appeared on its surface. The dragon’s eyes widened, and it let out a soft rumble of approval.
The egg was hatching!
Leo and Lily watched in amazement as the shell broke apart, and a tiny dragon emerged. It was no bigger than a rabbit, with shimmering silver scales and bright blue eyes. It looked up at them, blinking in the sunlight, and let out a
'@
$ws.Range("B5").Value = $b5Text
$ws.Rows.Item(5).AutoFit() | Out-Null

# B6: clear the cheese-making answer entirely (cell becomes empty)
$ws.Range("B6").ClearContents()

# Move the active selection from A12 to B12
$ws.Range("B12").Select() | Out-Null
